# Apply the "stats -> GET /api/stats" endpoint-spec rewrite to both sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "stats": drop the param:unitName column (I) and reword the test
# description now that the endpoint is described via HTTP method/path.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("stats")

$ws1.Range("B2").Value = "Test GET /api/stats endpoint with valid parameters"

# Remove column I entirely (header "param:unitName" and its blank value),
# which also shrinks the used range from A1:I2 down to A1:H2.
$ws1.Columns.Item(9).Delete()

# ---------------------------------------------------------------------
# Sheet "Documentation": swap the "API Function/File Path" framing for
# an "HTTP Method/Path/Description" framing, and drop the now-unused
# "Parameter Descriptions" block, shifting the remaining rows up.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Documentation")

$ws2.Range("A3").Value = "HTTP Method:"
$ws2.Range("B3").Value = "GET"

$ws2.Range("A4").Value = "Path:"
$ws2.Range("B4").Value = "/api/stats"

$ws2.Range("A5").Value = "Description:"
$ws2.Range("B5").Value = "Get system statistics"

$ws2.Range("A6").Value = ""

$ws2.Range("A7").Value = "Column Descriptions:"
$ws2.Range("B7").ClearContents()

$ws2.Range("A8").Value = "test_name"
$ws2.Range("B8").Value = "Unique identifier for the test"

$ws2.Range("A9").Value = "description"
$ws2.Range("B9").Value = "Human-readable description of what the test does"

$ws2.Range("A10").Value = "enabled"
$ws2.Range("B10").Value = "Whether to run this test (true/false)"

$ws2.Range("A11").Value = "expected_status"
$ws2.Range("B11").Value = "Expected HTTP status code (200, 404, etc.)"

$ws2.Range("A12").Value = "timeout_ms"
$ws2.Range("B12").Value = "Request timeout in milliseconds"

$ws2.Range("A13").Value = "max_response_time"
$ws2.Range("B13").Value = "Maximum acceptable response time in ms"

$ws2.Range("A14").Value = "delay_after_ms"
$ws2.Range("B14").Value = "Delay after test completion in ms"

$ws2.Range("A15").Value = "tags"
$ws2.Range("B15").Value = "Comma-separated tags for filtering tests"

$ws2.Range("A16").Value = ""

# Former "Parameter Descriptions:" / "param:unitName" rows collapse to a
# single blank row now that there are no more endpoint parameters.
$ws2.Range("A17").Value = ""
$ws2.Range("B17").ClearContents()

$ws2.Range("A18").Value = "Endpoint-Specific Notes:"
$ws2.Range("A19").Value = "• Method: GET"
$ws2.Range("A20").Value = "• Path: /api/stats"
$ws2.Range("A21").Value = "• Description: Get system statistics"
